$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 16.33600000000001
$ws.Range("C7").Value = -13.385
$ws.Range("B8").Value = 6.583000000000001
$ws.Range("B10").Value = 5.828999999999999
$ws.Range("B12").Value = 5.315
$ws.Range("C15").Value = -13.636
$ws.Range("B18").Value = 5.762
$ws.Range("C18").Value = -12.101
$ws.Range("E18").Value = 17.323
$ws.Range("E19").Value = 16.538
$ws.Range("C20").Value = -12.183
$ws.Range("E27").Value = 16.66
$ws.Range("C29").Value = -12.334
$ws.Range("C30").Value = -13.177
$ws.Range("C31").Value = -13.185
$ws.Range("E31").Value = 16.165
$ws.Range("B37").Value = 8.494
$ws.Range("E38").Value = 16.568
$ws.Range("C40").Value = -12.782
$ws.Range("E42").Value = 16.466
$ws.Range("E44").Value = 16.779
$ws.Range("E47").Value = 16.426
$ws.Range("C50").Value = -13.371
$ws.Range("B55").Value = 5.257
$ws.Range("E58").Value = 16.541
$ws.Range("E65").Value = 17.31
$ws.Range("B68").Value = 5.431
$ws.Range("C68").Value = -11.177
$ws.Range("E73").Value = 16.637
$ws.Range("C76").Value = -13.371
$ws.Range("B77").Value = 6.377
$ws.Range("B78").Value = 7.58
$ws.Range("B81").Value = 5.877
$ws.Range("B82").Value = 5.659
$ws.Range("C87").Value = -13.029
$ws.Range("C88").Value = -12.774
$ws.Range("E90").Value = 16.567
$ws.Range("E94").Value = 17.898
$ws.Range("E95").Value = 17.602
$ws.Range("C96").Value = -12.85
$ws.Range("C98").Value = -13.201
$ws.Range("C101").Value = -12.612
$ws.Range("E101").Value = 16.655
$ws.Range("C102").Value = -13.048
